$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data (columns F:V) between rows 36 and 37 ---
$ws.Range("F36").Value = "Youssoufia Berrechid"
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = "Mouloudia Oujda"
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 2.28
$ws.Range("K36").Value = "01/10/2023 04:42"
$ws.Range("L36").Value = 2.05
$ws.Range("M36").Value = "01/10/2023 19:11"
$ws.Range("N36").Value = 2.9
$ws.Range("O36").Value = "01/10/2023 04:42"
$ws.Range("P36").Value = 2.99
$ws.Range("Q36").Value = "01/10/2023 19:11"
$ws.Range("R36").Value = 3.29
$ws.Range("S36").Value = "01/10/2023 04:42"
$ws.Range("T36").Value = 4.02
$ws.Range("U36").Value = "01/10/2023 19:11"
$ws.Range("V36").Value = "https://www.betexplorer.com/football/morocco/botola-pro/youssoufia-berrechid-mouloudia-oujda/IeJ5XXs7/"

$ws.Range("F37").Value = "Maghreb Fez"
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = "Jeunesse Sportive Soualem"
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 2.09
$ws.Range("K37").Value = "30/09/2023 06:42"
$ws.Range("L37").Value = 2.01
$ws.Range("M37").Value = "01/10/2023 19:13"
$ws.Range("N37").Value = 2.82
$ws.Range("O37").Value = "30/09/2023 06:42"
$ws.Range("P37").Value = 3.03
$ws.Range("Q37").Value = "01/10/2023 19:01"
$ws.Range("R37").Value = 3.56
$ws.Range("S37").Value = "30/09/2023 06:42"
$ws.Range("T37").Value = 4.13
$ws.Range("U37").Value = "01/10/2023 19:13"
$ws.Range("V37").Value = "https://www.betexplorer.com/football/morocco/botola-pro/maghreb-fez-jeunesse-sportive-soualem/OYJ1YDS0/"

# --- Swap match data (columns F:V) between rows 47 and 48 ---
$ws.Range("F47").Value = "Renaissance Zemamra"
$ws.Range("G47").Value = 3
$ws.Range("H47").Value = "Olympique de Safi"
$ws.Range("I47").Value = 2
$ws.Range("J47").Value = 2.98
$ws.Range("K47").Value = "08/10/2023 04:12"
$ws.Range("L47").Value = 3.95
$ws.Range("M47").Value = "08/10/2023 19:14"
$ws.Range("N47").Value = 2.85
$ws.Range("O47").Value = "08/10/2023 04:12"
$ws.Range("P47").Value = 2.75
$ws.Range("Q47").Value = "08/10/2023 19:14"
$ws.Range("R47").Value = 2.5
$ws.Range("S47").Value = "08/10/2023 04:12"
$ws.Range("T47").Value = 2.2
$ws.Range("U47").Value = "08/10/2023 19:14"
$ws.Range("V47").Value = "https://www.betexplorer.com/football/morocco/botola-pro/renaissance-zemamra-olympique-de-safi/xSSTnVt3/"

$ws.Range("F48").Value = "Berkane"
$ws.Range("G48").Value = 3
$ws.Range("H48").Value = "Youssoufia Berrechid"
$ws.Range("I48").Value = 1
$ws.Range("J48").Value = 1.39
$ws.Range("K48").Value = "08/10/2023 04:12"
$ws.Range("L48").Value = 1.33
$ws.Range("M48").Value = "08/10/2023 18:54"
$ws.Range("N48").Value = 4.21
$ws.Range("O48").Value = "08/10/2023 04:12"
$ws.Range("P48").Value = 4.64
$ws.Range("Q48").Value = "08/10/2023 18:54"
$ws.Range("R48").Value = 7.84
$ws.Range("S48").Value = "08/10/2023 04:12"
$ws.Range("T48").Value = 9.8
$ws.Range("U48").Value = "08/10/2023 18:54"
$ws.Range("V48").Value = "https://www.betexplorer.com/football/morocco/botola-pro/berkane-youssoufia-berrechid/W2sXokd9/"

# --- Append new row 76 (new match: FAR Rabat vs Wydad) ---
$ws.Range("A75:V75").Copy()
$ws.Range("A76:V76").PasteSpecial(-4122)
$ws.Range("A76").Value = 75
$ws.Range("B76").Value = "morocco"
$ws.Range("C76").Value = "botola-pro"
$ws.Range("D76").Value = "2023-2024"
$ws.Range("E76").Value = 45258.75
$ws.Range("F76").Value = "FAR Rabat"
$ws.Range("G76").Value = 3
$ws.Range("H76").Value = "Wydad"
$ws.Range("I76").Value = 1
$ws.Range("J76").Value = 2.18
$ws.Range("K76").Value = "26/11/2023 16:13"
$ws.Range("L76").Value = 2.29
$ws.Range("M76").Value = "28/11/2023 17:51"
$ws.Range("N76").Value = 2.88
$ws.Range("O76").Value = "26/11/2023 16:13"
$ws.Range("P76").Value = 2.8
$ws.Range("Q76").Value = "28/11/2023 17:51"
$ws.Range("R76").Value = 3.24
$ws.Range("S76").Value = "26/11/2023 16:13"
$ws.Range("T76").Value = 3.61
$ws.Range("U76").Value = "28/11/2023 17:51"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/morocco/botola-pro/far-rabat-wydad-athletic/lMW0H6E6/"
